# Update the multiplication "fact family" answers in the practice-sheet
# table. Every left-hand value below is unique in the document, so a
# literal (non-wildcard) Find/Replace targets exactly one cell each,
# leaving all other formatting (fonts, sizes, table layout) untouched.
#
# Find.Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,
#              MatchSoundsLike, MatchAllWordForms, Forward, Wrap,
#              Format, ReplaceWith, Replace)
#   Wrap = 1 (wdFindContinue), Replace = 2 (wdReplaceAll)

$d = $word.ActiveDocument

$d.Content.Find.Execute("69×39=2691", $true, $false, $false, $false, $false, $true, 1, $false, "67×97=6499", 2)
$d.Content.Find.Execute("88×90=7920", $true, $false, $false, $false, $false, $true, 1, $false, "42×17=714", 2)
$d.Content.Find.Execute("64×14=896", $true, $false, $false, $false, $false, $true, 1, $false, "66×39=2574", 2)
$d.Content.Find.Execute("84×22=1848", $true, $false, $false, $false, $false, $true, 1, $false, "17×55=935", 2)
$d.Content.Find.Execute("23×66=1518", $true, $false, $false, $false, $false, $true, 1, $false, "22×13=286", 2)
$d.Content.Find.Execute("97×20=1940", $true, $false, $false, $false, $false, $true, 1, $false, "88×46=4048", 2)
$d.Content.Find.Execute("93×97=9021", $true, $false, $false, $false, $false, $true, 1, $false, "29×81=2349", 2)
$d.Content.Find.Execute("25×73=1825", $true, $false, $false, $false, $false, $true, 1, $false, "97×82=7954", 2)
$d.Content.Find.Execute("47×55=2585", $true, $false, $false, $false, $false, $true, 1, $false, "99×23=2277", 2)
$d.Content.Find.Execute("32×30=960", $true, $false, $false, $false, $false, $true, 1, $false, "37×74=2738", 2)
$d.Content.Find.Execute("47×50=2350", $true, $false, $false, $false, $false, $true, 1, $false, "51×39=1989", 2)
$d.Content.Find.Execute("85×12=1020", $true, $false, $false, $false, $false, $true, 1, $false, "98×59=5782", 2)
$d.Content.Find.Execute("22×14=308", $true, $false, $false, $false, $false, $true, 1, $false, "41×70=2870", 2)
$d.Content.Find.Execute("44×55=2420", $true, $false, $false, $false, $false, $true, 1, $false, "24×62=1488", 2)
$d.Content.Find.Execute("27×18=486", $true, $false, $false, $false, $false, $true, 1, $false, "33×65=2145", 2)
$d.Content.Find.Execute("22×49=1078", $true, $false, $false, $false, $false, $true, 1, $false, "12×68=816", 2)
$d.Content.Find.Execute("39×41=1599", $true, $false, $false, $false, $false, $true, 1, $false, "99×33=3267", 2)
$d.Content.Find.Execute("80×13=1040", $true, $false, $false, $false, $false, $true, 1, $false, "23×69=1587", 2)
$d.Content.Find.Execute("82×99=8118", $true, $false, $false, $false, $false, $true, 1, $false, "16×73=1168", 2)
$d.Content.Find.Execute("57×68=3876", $true, $false, $false, $false, $false, $true, 1, $false, "52×96=4992", 2)
$d.Content.Find.Execute("57×61=3477", $true, $false, $false, $false, $false, $true, 1, $false, "63×99=6237", 2)
$d.Content.Find.Execute("54×76=4104", $true, $false, $false, $false, $false, $true, 1, $false, "67×96=6432", 2)
$d.Content.Find.Execute("64×24=1536", $true, $false, $false, $false, $false, $true, 1, $false, "31×75=2325", 2)
$d.Content.Find.Execute("77×20=1540", $true, $false, $false, $false, $false, $true, 1, $false, "51×29=1479", 2)
$d.Content.Find.Execute("22×41=902", $true, $false, $false, $false, $false, $true, 1, $false, "56×53=2968", 2)

$d.Save()

